$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1156.8572
$ws.Range("J33").Value = 3091
$ws.Range("L33").Value = 3091
$ws.Range("N33").Value = -3549
$ws.Range("H61").Value = 4188130
$ws.Range("I61").Value = 4762146
$ws.Range("J61").Value = 170017
$ws.Range("K61").Value = 14286438
$ws.Range("L61").Value = 510051
$ws.Range("M61").Value = -14286266
$ws.Range("N61").Value = -510395
$ws.Range("H62").Value = 1916.2778
$ws.Range("I62").Value = 1974.5834
$ws.Range("K62").Value = 1974.5834
$ws.Range("M62").Value = -1350.5834
$ws.Range("H65").Value = 1916.2778
$ws.Range("I65").Value = 1974.5834
$ws.Range("K65").Value = 9872.916999999999
$ws.Range("M65").Value = -6752.916999999999
$ws.Range("H100").Value = 1372.5714
$ws.Range("I100").Value = 675
$ws.Range("J100").Value = 1895.75
$ws.Range("K100").Value = 675
$ws.Range("L100").Value = 1895.75
$ws.Range("M100").Value = -134
$ws.Range("N100").Value = -2977.75
$ws.Range("H112").Value = 1179.3793
$ws.Range("J112").Value = 1207.909
$ws.Range("L112").Value = 3623.727
$ws.Range("N112").Value = -5839.727000000001
$ws.Range("H132").Value = 5561401
$ws.Range("I132").Value = 6416155
$ws.Range("K132").Value = 19248465
$ws.Range("M132").Value = -19245935
$ws.Range("H138").Value = 2731.3333
$ws.Range("I138").Value = 1424.6207
$ws.Range("J138").Value = 3755.5134
$ws.Range("K138").Value = 4273.8621
$ws.Range("L138").Value = 11266.5402
$ws.Range("M138").Value = 866.1378999999997
$ws.Range("N138").Value = -21546.5402

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 145351.14
$ws.Range("I102").Value = 252120
$ws.Range("J102").Value = 2992.6667
$ws.Range("K102").Value = 252120
$ws.Range("L102").Value = 2992.6667
$ws.Range("M102").Value = -250498
$ws.Range("N102").Value = -6236.6667
$ws.Range("H110").Value = 167017340
$ws.Range("I110").Value = 167017340
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 167017340
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -167015295
$ws.Range("N110").Value = ""
$ws.Range("H122").Value = 2374.625
$ws.Range("I122").Value = 2001.9412
$ws.Range("K122").Value = 6005.8236
$ws.Range("M122").Value = -3555.8236
$ws.Range("H132").Value = 16192.244
$ws.Range("I132").Value = 19702.656
$ws.Range("J132").Value = 3710.7778
$ws.Range("K132").Value = 59107.96799999999
$ws.Range("L132").Value = 11132.3334
$ws.Range("M132").Value = -56577.96799999999
$ws.Range("N132").Value = -16192.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 51709.773
$ws.Range("I86").Value = 62772.668
$ws.Range("J86").Value = 1926.75
$ws.Range("K86").Value = 62772.668
$ws.Range("L86").Value = 1926.75
$ws.Range("M86").Value = -61649.668
$ws.Range("N86").Value = -4172.75
$ws.Range("H89").Value = 51709.773
$ws.Range("I89").Value = 62772.668
$ws.Range("J89").Value = 1926.75
$ws.Range("K89").Value = 313863.34
$ws.Range("L89").Value = 9633.75
$ws.Range("M89").Value = -308247.34
$ws.Range("N89").Value = -20865.75
$ws.Range("H107").Value = 55556556
$ws.Range("I107").Value = 66667650
$ws.Range("J107").Value = 1096.6666
$ws.Range("K107").Value = 66667650
$ws.Range("L107").Value = 1096.6666
$ws.Range("M107").Value = -66665730
$ws.Range("N107").Value = -4936.6666
$ws.Range("H134").Value = 1631.963
$ws.Range("I134").Value = 1311.94
$ws.Range("J134").Value = 5632.25
$ws.Range("K134").Value = 3935.82
$ws.Range("L134").Value = 16896.75
$ws.Range("M134").Value = -1400.82
$ws.Range("N134").Value = -21966.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 1518.75
$ws.Range("J8").Value = 1758.3334
$ws.Range("L8").Value = 1758.3334
$ws.Range("N8").Value = -2038.3334
$ws.Range("H25").Value = 20000
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = ""
$ws.Range("H31").Value = 37384.023
$ws.Range("J31").Value = 44898.383
$ws.Range("L31").Value = 44898.383
$ws.Range("N31").Value = -45488.383
$ws.Range("H34").Value = 37384.023
$ws.Range("J34").Value = 44898.383
$ws.Range("L34").Value = 44898.383
$ws.Range("N34").Value = -45302.383
$ws.Range("H62").Value = 2609
$ws.Range("I62").Value = 2450
$ws.Range("J62").Value = 2644.3333
$ws.Range("K62").Value = 2450
$ws.Range("L62").Value = 2644.3333
$ws.Range("M62").Value = -1826
$ws.Range("N62").Value = -3892.3333
$ws.Range("H65").Value = 2609
$ws.Range("I65").Value = 2450
$ws.Range("J65").Value = 2644.3333
$ws.Range("K65").Value = 12250
$ws.Range("L65").Value = 13221.6665
$ws.Range("M65").Value = -9130
$ws.Range("N65").Value = -19461.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 681
$ws.Range("J2").Value = 1012.3
$ws.Range("L2").Value = 6073.799999999999
$ws.Range("N2").Value = -6299.799999999999
$ws.Range("H3").Value = 992.7273
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = ""
$ws.Range("H5").Value = 1399.6346
$ws.Range("J5").Value = 1339.75
$ws.Range("L5").Value = 4019.25
$ws.Range("N5").Value = -4243.25
$ws.Range("H15").Value = 95.80952499999999
$ws.Range("I15").Value = 31.333334
$ws.Range("J15").Value = 257
$ws.Range("K15").Value = 94.00000199999999
$ws.Range("L15").Value = 771
$ws.Range("M15").Value = 45.99999800000001
$ws.Range("N15").Value = -1051
$ws.Range("H20").Value = 900
$ws.Range("I20").Value = 900
$ws.Range("K20").Value = 2700
$ws.Range("M20").Value = -2473
$ws.Range("H22").Value = 13561.5
$ws.Range("J22").Value = 13561.5
$ws.Range("L22").Value = 40684.5
$ws.Range("N22").Value = -41022.5
$ws.Range("H27").Value = 13561.5
$ws.Range("J27").Value = 13561.5
$ws.Range("L27").Value = 40684.5
$ws.Range("N27").Value = -40888.5
$ws.Range("H37").Value = 5012003
$ws.Range("J37").Value = 5012003
$ws.Range("L37").Value = 15036009
$ws.Range("N37").Value = -15036233
$ws.Range("H100").Value = 2500
$ws.Range("J100").Value = 2500
$ws.Range("L100").Value = 7500
$ws.Range("N100").Value = -9122
$ws.Range("H129").Value = 5880.696
$ws.Range("I129").Value = 622.25
$ws.Range("J129").Value = 6987.737
$ws.Range("K129").Value = 1866.75
$ws.Range("L129").Value = 20963.211
$ws.Range("M129").Value = 3133.25
$ws.Range("N129").Value = -30963.211
$ws.Range("H135").Value = 1399.6346
$ws.Range("J135").Value = 1339.75
$ws.Range("L135").Value = 12057.75
$ws.Range("N135").Value = -17127.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4086.7917
$ws.Range("I132").Value = 3098
$ws.Range("J132").Value = 5471.1
$ws.Range("K132").Value = 9294
$ws.Range("L132").Value = 16413.3
$ws.Range("M132").Value = -6764
$ws.Range("N132").Value = -21473.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3350.423
$ws.Range("I132").Value = 3444.44
$ws.Range("K132").Value = 10333.32
$ws.Range("M132").Value = -7803.32
$ws.Range("H140").Value = 99114.5
$ws.Range("J140").Value = 99114.5
$ws.Range("L140").Value = 99114.5
$ws.Range("N140").Value = -109474.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 8661.5
$ws.Range("I45").Value = 5000
$ws.Range("J45").Value = 9184.571
$ws.Range("K45").Value = 5000
$ws.Range("L45").Value = 9184.571
$ws.Range("M45").Value = -4509
$ws.Range("N45").Value = -10166.571
$ws.Range("H74").Value = 8305.857
$ws.Range("J74").Value = 8305.857
$ws.Range("L74").Value = 8305.857
$ws.Range("N74").Value = -10177.857
$ws.Range("H77").Value = 8305.857
$ws.Range("J77").Value = 8305.857
$ws.Range("L77").Value = 24917.571
$ws.Range("N77").Value = -34277.571
$ws.Range("H132").Value = 9493.700000000001
$ws.Range("I132").Value = 5402.769
$ws.Range("J132").Value = 17091.143
$ws.Range("K132").Value = 16208.307
$ws.Range("L132").Value = 51273.429
$ws.Range("M132").Value = -13678.307
$ws.Range("N132").Value = -56333.429
$ws.Range("H136").Value = 20316.29
$ws.Range("I136").Value = 39144.08
$ws.Range("J136").Value = 5482.273
$ws.Range("K136").Value = 117432.24
$ws.Range("L136").Value = 16446.819
$ws.Range("M136").Value = -114882.24
$ws.Range("N136").Value = -21546.819
